$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the trailing comment text first (while it is still on its
# original row, A16) so the existing shared string is edited in place
# instead of a new one being appended after the new "invalid"/"valid"/
# "players = null" strings introduced below.
$ws.Range("A16").Value = "// Vi kommer behöva förklara att inga parametrar är en valid ekvivalensklass också"

# --- Table 1 (subtractChips): insert a new "invalid/valid/invalid" header
# row above the existing "Testfalls ID" row, pushing everything below down
# by one row.
$ws.Rows("2:2").Insert()
$ws.Range("B2").Value = "invalid"
$ws.Range("C2").Value = "valid"
$ws.Range("D2").Value = "invalid"

# --- Table 2 (game): insert a new "invalid/valid/.." header row above the
# "Testfalls ID" row of the second table (originally row 9, now row 10
# after the first insert above), pushing the rest of that table down by
# one row.
$ws.Rows("10:10").Insert()

# Insert a new column (G) for the "players = null" equivalence class, which
# pushes the old "Valid?" column (and its data) from G to H.
$ws.Columns("G:G").Insert()

# The new column should carry the same custom width as column F (16).
$ws.Range("G1").ColumnWidth = $ws.Range("F1").ColumnWidth

# Now fill in the new header row for table 2 (columns B-G).
$ws.Range("B10").Value = "invalid"
$ws.Range("C10").Value = "valid"
$ws.Range("D10").Value = "invalid"
$ws.Range("E10").Value = "valid"
$ws.Range("F10").Value = "valid"
$ws.Range("G10").Value = "invalid"

# Label the newly inserted column.
$ws.Range("G11").Value = "players = null"

# The "x" marks that used to live in the old column G (now shifted to
# column H by the insert) for test cases 4 and 5 no longer apply, so clear
# them back out.
$ws.Range("H15").ClearContents()
$ws.Range("H16").ClearContents()

# Add the new 6th test case row for table 2 into what was the blank
# "gap" row separating the table from the trailing comment (now row 17,
# after the earlier inserts) - no extra row needs to be inserted, the
# comment row naturally stays put right below it at row 18.
$ws.Range("A17").Value = 6
$ws.Range("C17").Value = "x"
$ws.Range("E17").Value = "x"
$ws.Range("G17").Value = "x"

# Restore the selection to match the saved view state.
$ws.Range("H16").Select()
